$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.096431111944355052,
    -0.0099999993200050596,
    -0.067292533817552069,
    0.28398800098384314,
    -0.0059999993475994273,
    -0.0059999993239294724,
    -0.019999999192407358,
    -0.019999999182501504,
    -0.0059999992986172757,
    -0.0059999992872548091,
    -0.0044999992996253013,
    -0.0059999992832167059,
    -0.005999999271095291,
    -0.011999999212460644,
    -0.0059999992645982658,
    -0.0059999992620165532,
    -0.0059999992585604289,
    -0.0089999992303013698,
    -0.0089999993291827174,
    -0.0089999993231462128,
    -0.0089999993222455998,
    -0.0089999993215839069,
    -0.0089999993185898575,
    -0.041999999005571276,
    -0.041999999000113419,
    -0.0059999993229880033,
    -0.0059999993215633651,
    -0.0059999993154882247,
    -0.01199999925692552,
    -0.019999999182596984,
    -0.014447580125992232,
    -0.020999999174823536,
    0.039310346265988727
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
